$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5747
$ws.Range("I40").Value = 9850
$ws.Range("J40").Value = 4721.25
$ws.Range("K40").Value = 9850
$ws.Range("L40").Value = 4721.25
$ws.Range("M40").Value = -9675
$ws.Range("N40").Value = -5071.25
$ws.Range("H64").Value = 6253
$ws.Range("J64").Value = 11000
$ws.Range("L64").Value = 11000
$ws.Range("N64").Value = -11496
$ws.Range("H67").Value = 6253
$ws.Range("J67").Value = 11000
$ws.Range("L67").Value = 11000
$ws.Range("N67").Value = -12716
$ws.Range("H74").Value = 4228.6665
$ws.Range("I74").Value = 4228.6665
$ws.Range("K74").Value = 4228.6665
$ws.Range("M74").Value = -3292.6665
$ws.Range("H77").Value = 4228.6665
$ws.Range("I77").Value = 4228.6665
$ws.Range("K77").Value = 21143.3325
$ws.Range("M77").Value = -16463.3325
$ws.Range("H98").Value = 1289.3158
$ws.Range("I98").Value = 1334.3334
$ws.Range("J98").Value = 479
$ws.Range("K98").Value = 1334.3334
$ws.Range("L98").Value = 479
$ws.Range("M98").Value = 163.6666
$ws.Range("N98").Value = -3475
$ws.Range("H122").Value = 1289.3158
$ws.Range("I122").Value = 1334.3334
$ws.Range("J122").Value = 479
$ws.Range("K122").Value = 4003.0002
$ws.Range("L122").Value = 1437
$ws.Range("M122").Value = -1553.0002
$ws.Range("N122").Value = -6337
$ws.Range("H137").Value = 667.55554
$ws.Range("I137").Value = 667.55554
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2002.66662
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 547.33338
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 1679.4706
$ws.Range("I138").Value = 1530.7333
$ws.Range("J138").Value = 2795
$ws.Range("K138").Value = 4592.199900000001
$ws.Range("L138").Value = 8385
$ws.Range("M138").Value = 547.8000999999995
$ws.Range("N138").Value = -18665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3915.9714
$ws.Range("I32").Value = 3915.9714
$ws.Range("K32").Value = 3915.9714
$ws.Range("M32").Value = -3628.9714
$ws.Range("H63").Value = 5841.4287
$ws.Range("I63").Value = 4478
$ws.Range("J63").Value = 9250
$ws.Range("K63").Value = 4478
$ws.Range("L63").Value = 9250
$ws.Range("M63").Value = -3792
$ws.Range("N63").Value = -10622
$ws.Range("H66").Value = 5841.4287
$ws.Range("I66").Value = 4478
$ws.Range("J66").Value = 9250
$ws.Range("K66").Value = 22390
$ws.Range("L66").Value = 46250
$ws.Range("M66").Value = -18958
$ws.Range("N66").Value = -53114
$ws.Range("H74").Value = 1347.5
$ws.Range("I74").Value = 700
$ws.Range("J74").Value = 1995
$ws.Range("K74").Value = 700
$ws.Range("L74").Value = 1995
$ws.Range("M74").Value = 174
$ws.Range("N74").Value = -3743
$ws.Range("H77").Value = 1347.5
$ws.Range("I77").Value = 700
$ws.Range("J77").Value = 1995
$ws.Range("K77").Value = 3500
$ws.Range("L77").Value = 9975
$ws.Range("M77").Value = 868
$ws.Range("N77").Value = -18711

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4125.8237
$ws.Range("I86").Value = 2995.7
$ws.Range("J86").Value = 5740.2856
$ws.Range("K86").Value = 2995.7
$ws.Range("L86").Value = 5740.2856
$ws.Range("M86").Value = -1872.7
$ws.Range("N86").Value = -7986.2856
$ws.Range("H89").Value = 4125.8237
$ws.Range("I89").Value = 2995.7
$ws.Range("J89").Value = 5740.2856
$ws.Range("K89").Value = 14978.5
$ws.Range("L89").Value = 28701.428
$ws.Range("M89").Value = -9362.5
$ws.Range("N89").Value = -39933.428
$ws.Range("H105").Value = 2701.1428
$ws.Range("I105").Value = 2681.6
$ws.Range("K105").Value = 2681.6
$ws.Range("M105").Value = -934.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3973.3076
$ws.Range("I31").Value = 2064.25
$ws.Range("J31").Value = 4821.778
$ws.Range("K31").Value = 2064.25
$ws.Range("L31").Value = 4821.778
$ws.Range("M31").Value = -1769.25
$ws.Range("N31").Value = -5411.778
$ws.Range("H34").Value = 3973.3076
$ws.Range("I34").Value = 2064.25
$ws.Range("J34").Value = 4821.778
$ws.Range("K34").Value = 2064.25
$ws.Range("L34").Value = 4821.778
$ws.Range("M34").Value = -1862.25
$ws.Range("N34").Value = -5225.778
$ws.Range("H62").Value = 2439.2
$ws.Range("I62").Value = 2065.3333
$ws.Range("K62").Value = 2065.3333
$ws.Range("M62").Value = -1441.3333
$ws.Range("H65").Value = 2439.2
$ws.Range("I65").Value = 2065.3333
$ws.Range("K65").Value = 10326.6665
$ws.Range("M65").Value = -7206.666499999999
$ws.Range("H94").Value = 1
$ws.Range("J94").Value = 1
$ws.Range("L94").Value = 1
$ws.Range("N94").Value = -903
$ws.Range("H107").Value = 946.125
$ws.Range("I107").Value = 595
$ws.Range("K107").Value = 595
$ws.Range("M107").Value = 1325
$ws.Range("H141").Value = 531250
$ws.Range("J141").Value = 1000000
$ws.Range("L141").Value = 1000000
$ws.Range("N141").Value = -1010360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 58.166668
$ws.Range("J38").Value = 20
$ws.Range("L38").Value = 60
$ws.Range("N38").Value = -754
$ws.Range("H80").Value = 8999.25
$ws.Range("I80").Value = 8999
$ws.Range("J80").Value = 8999.333000000001
$ws.Range("K80").Value = 26997
$ws.Range("L80").Value = 26997.999
$ws.Range("M80").Value = -26061
$ws.Range("N80").Value = -28869.999
$ws.Range("H83").Value = 8999.25
$ws.Range("I83").Value = 8999
$ws.Range("J83").Value = 8999.333000000001
$ws.Range("K83").Value = 80991
$ws.Range("L83").Value = 80993.997
$ws.Range("M83").Value = -76311
$ws.Range("N83").Value = -90353.997
$ws.Range("H114").Value = 5613.5
$ws.Range("I114").Value = 5227.5
$ws.Range("J114").Value = 5999.5
$ws.Range("K114").Value = 15682.5
$ws.Range("L114").Value = 17998.5
$ws.Range("M114").Value = -12428.5
$ws.Range("N114").Value = -24506.5
$ws.Range("H140").Value = 912779
$ws.Range("I140").Value = 912779
$ws.Range("K140").Value = 2738337
$ws.Range("M140").Value = -2733157

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3407.8572
$ws.Range("J80").Value = 4350
$ws.Range("L80").Value = 4350
$ws.Range("N80").Value = -6346
$ws.Range("H83").Value = 3407.8572
$ws.Range("J83").Value = 4350
$ws.Range("L83").Value = 21750
$ws.Range("N83").Value = -31734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 588888.9
$ws.Range("J43").Value = 588888.9
$ws.Range("L43").Value = 588888.9
$ws.Range("N43").Value = -589274.9
$ws.Range("H55").Value = 666.75
$ws.Range("I55").Value = 472.33334
$ws.Range("K55").Value = 472.33334
$ws.Range("M55").Value = -299.33334
$ws.Range("H82").Value = 2344.6667
$ws.Range("I82").Value = 1200
$ws.Range("J82").Value = 2917
$ws.Range("K82").Value = 1200
$ws.Range("L82").Value = 2917
$ws.Range("M82").Value = -839
$ws.Range("N82").Value = -3639
$ws.Range("H85").Value = 2344.6667
$ws.Range("I85").Value = 1200
$ws.Range("J85").Value = 2917
$ws.Range("K85").Value = 1200
$ws.Range("L85").Value = 2917
$ws.Range("M85").Value = 48
$ws.Range("N85").Value = -5413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3770875
$ws.Range("J4").Value = 3770875
$ws.Range("L4").Value = 3770875
$ws.Range("N4").Value = -3771101
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H124").Value = 23946.666
$ws.Range("J124").Value = 23946.666
$ws.Range("L124").Value = 23946.666
$ws.Range("N124").Value = -33766.666
$ws.Range("H132").Value = 4740
$ws.Range("I132").Value = 3425
$ws.Range("K132").Value = 10275
$ws.Range("M132").Value = -7745
$ws.Range("H136").Value = 19029.8
$ws.Range("I136").Value = 16271.286
$ws.Range("J136").Value = 25466.334
$ws.Range("K136").Value = 48813.858
$ws.Range("L136").Value = 76399.00199999999
$ws.Range("M136").Value = -46263.858
$ws.Range("N136").Value = -81499.00199999999

Write-Host "Applied all changes"